# Auto-generated script applying the XML diff changes to Ifrit_Profits workbook
# Updates numeric values for columns H-N on specific rows across all 8 job sheets,
# including removing obsolete LeveProfitHQ (N) cells on LTW rows 68/71/100/132,
# and adding a new LeveProfitHQ (N139) cell on the CUL sheet.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 199.1
$ws.Range("I38").Value = 110.111115
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 330.333345
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = 41.66665499999999
$ws.Range("N38").Value = -3744
$ws.Range("H40").Value = 950.35
$ws.Range("I40").Value = 941.64703
$ws.Range("J40").Value = 999.6667
$ws.Range("K40").Value = 941.64703
$ws.Range("L40").Value = 999.6667
$ws.Range("M40").Value = -766.64703
$ws.Range("N40").Value = -1349.6667
$ws.Range("H112").Value = 23810726
$ws.Range("J112").Value = 24391462
$ws.Range("L112").Value = 73174386
$ws.Range("N112").Value = -73176602
$ws.Range("H116").Value = 2499.5454
$ws.Range("I116").Value = 2142.8572
$ws.Range("J116").Value = 2666
$ws.Range("K116").Value = 2142.8572
$ws.Range("L116").Value = 2666
$ws.Range("M116").Value = 1299.1428
$ws.Range("N116").Value = -9550

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 155365.08
$ws.Range("I45").Value = 168212.17
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 168212.17
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -167835.17
$ws.Range("N45").Value = -1954
$ws.Range("H132").Value = 1718.8837
$ws.Range("I132").Value = 1592.2222
$ws.Range("J132").Value = 1932.625
$ws.Range("K132").Value = 4776.6666
$ws.Range("L132").Value = 5797.875
$ws.Range("M132").Value = -2246.6666
$ws.Range("N132").Value = -10857.875

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2303.5386
$ws.Range("I86").Value = 2226.5334
$ws.Range("J86").Value = 2408.5454
$ws.Range("K86").Value = 2226.5334
$ws.Range("L86").Value = 2408.5454
$ws.Range("M86").Value = -1103.5334
$ws.Range("N86").Value = -4654.5454
$ws.Range("H89").Value = 2303.5386
$ws.Range("I89").Value = 2226.5334
$ws.Range("J89").Value = 2408.5454
$ws.Range("K89").Value = 11132.667
$ws.Range("L89").Value = 12042.727
$ws.Range("M89").Value = -5516.666999999999
$ws.Range("N89").Value = -23274.727
$ws.Range("H131").Value = 44056.668
$ws.Range("J131").Value = 44056.668
$ws.Range("L131").Value = 44056.668
$ws.Range("N131").Value = -54136.668

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2742.9119
$ws.Range("I31").Value = 1558.3846
$ws.Range("J31").Value = 3476.1904
$ws.Range("K31").Value = 1558.3846
$ws.Range("L31").Value = 3476.1904
$ws.Range("M31").Value = -1263.3846
$ws.Range("N31").Value = -4066.1904
$ws.Range("H34").Value = 2742.9119
$ws.Range("I34").Value = 1558.3846
$ws.Range("J34").Value = 3476.1904
$ws.Range("K34").Value = 1558.3846
$ws.Range("L34").Value = 3476.1904
$ws.Range("M34").Value = -1356.3846
$ws.Range("N34").Value = -3880.1904
$ws.Range("H132").Value = 2333.182
$ws.Range("I132").Value = 1491.1578
$ws.Range("J132").Value = 7666
$ws.Range("K132").Value = 4473.4734
$ws.Range("L132").Value = 22998
$ws.Range("M132").Value = -1943.4734
$ws.Range("N132").Value = -28058

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1254.3529
$ws.Range("I5").Value = 1338.1428
$ws.Range("J5").Value = 863.3333
$ws.Range("K5").Value = 4014.4284
$ws.Range("L5").Value = 2589.9999
$ws.Range("M5").Value = -3902.4284
$ws.Range("N5").Value = -2813.9999
$ws.Range("H117").Value = 2490
$ws.Range("I117").Value = 2000
$ws.Range("K117").Value = 6000
$ws.Range("M117").Value = -2558
$ws.Range("H129").Value = 886.2727
$ws.Range("I129").Value = 494.83334
$ws.Range("K129").Value = 1484.50002
$ws.Range("M129").Value = 3515.49998
$ws.Range("H131").Value = 1268039.6
$ws.Range("I131").Value = 5333.5
$ws.Range("J131").Value = 1451040.5
$ws.Range("K131").Value = 16000.5
$ws.Range("L131").Value = 4353121.5
$ws.Range("M131").Value = -10960.5
$ws.Range("N131").Value = -4363201.5
$ws.Range("H135").Value = 1254.3529
$ws.Range("I135").Value = 1338.1428
$ws.Range("J135").Value = 863.3333
$ws.Range("K135").Value = 12043.2852
$ws.Range("L135").Value = 7769.9997
$ws.Range("M135").Value = -9508.2852
$ws.Range("N135").Value = -12839.9997
$ws.Range("H139").Value = 2144.0557
$ws.Range("I139").Value = 1173.3334
$ws.Range("J139").Value = 3114.7778
$ws.Range("K139").Value = 3520.0002
$ws.Range("L139").Value = 9344.3334
$ws.Range("M139").Value = 1619.9998
$ws.Range("N139").Value = -19624.3334

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 241051.2
$ws.Range("I80").Value = 2900
$ws.Range("J80").Value = 300589
$ws.Range("K80").Value = 2900
$ws.Range("L80").Value = 300589
$ws.Range("M80").Value = -1902
$ws.Range("N80").Value = -302585
$ws.Range("H83").Value = 241051.2
$ws.Range("I83").Value = 2900
$ws.Range("J83").Value = 300589
$ws.Range("K83").Value = 14500
$ws.Range("L83").Value = 1502945
$ws.Range("M83").Value = -9508
$ws.Range("N83").Value = -1512929
$ws.Range("H132").Value = 2939.963
$ws.Range("I132").Value = 2800.2
$ws.Range("J132").Value = 3022.1765
$ws.Range("K132").Value = 8400.599999999999
$ws.Range("L132").Value = 9066.529500000001
$ws.Range("M132").Value = -5870.599999999999
$ws.Range("N132").Value = -14126.5295

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 925.125
$ws.Range("I46").Value = 962.75
$ws.Range("J46").Value = 887.5
$ws.Range("K46").Value = 962.75
$ws.Range("L46").Value = 887.5
$ws.Range("M46").Value = -774.75
$ws.Range("N46").Value = -1263.5
$ws.Range("H68").Value = 2033.3334
$ws.Range("I68").Value = 2033.3334
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2033.3334
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1284.3334
$ws.Range("H71").Value = 2033.3334
$ws.Range("I71").Value = 2033.3334
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 10166.667
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -6422.666999999999
$ws.Range("H82").Value = 1879.0769
$ws.Range("I82").Value = 1734.2858
$ws.Range("J82").Value = 2048
$ws.Range("K82").Value = 1734.2858
$ws.Range("L82").Value = 2048
$ws.Range("M82").Value = -1373.2858
$ws.Range("N82").Value = -2770
$ws.Range("H85").Value = 1879.0769
$ws.Range("I85").Value = 1734.2858
$ws.Range("J85").Value = 2048
$ws.Range("K85").Value = 1734.2858
$ws.Range("L85").Value = 2048
$ws.Range("M85").Value = -486.2858000000001
$ws.Range("N85").Value = -4544
$ws.Range("H100").Value = 1003
$ws.Range("I100").Value = 1003
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1003
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -462
$ws.Range("H132").Value = 189084.67
$ws.Range("I132").Value = 189084.67
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 567254.01
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -564724.01
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("N132").ClearContents()

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1393.4
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 1655.6666
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 1655.6666
$ws.Range("M96").Value = 373
$ws.Range("N96").Value = -4401.6666
$ws.Range("H136").Value = 965.67645
$ws.Range("I136").Value = 942.73914
$ws.Range("J136").Value = 1013.63635
$ws.Range("K136").Value = 2828.21742
$ws.Range("L136").Value = 3040.90905
$ws.Range("M136").Value = -278.2174199999999
$ws.Range("N136").Value = -8140.90905
